# 4COM02_ProteomicsDataProcessing.xlsx
# Add a "Parameter [analysis software]" block (with its Term Source REF /
# Term Accession Number columns) to the main annotation table, fill in
# some content examples for the acquisition/analysis/data-processing
# software columns, and bump the template version from 1.1.3 to 1.1.4.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Main annotation table (sheet 1 / "4COM02_ProteomicsDataProcessing")
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(1)
$lo = $ws.ListObjects.Item(1)

# Grow the table from 8 columns (A:H) / 2 rows to 11 columns (A:K) / 7 rows.
# New blank columns are appended at the tail; we rename/refill the cells
# below so the final layout matches the target (acquisition software,
# analysis software [new], data processing software, data file name).
$lo.Resize($ws.Range("A1:K7"))

# Re-label the header row into its final order:
# A Source Name | B Parameter [acquisition software]
# C Term Source REF (MS:1001455) | D Term Accession Number (MS:1001455)
# E Parameter [analysis software] (NEW)
# F Term Source REF (MS:1001456) (NEW) | G Term Accession Number (MS:1001456) (NEW)
# H Parameter [data processing software]
# I Term Source REF (MS:1001457) | J Term Accession Number (MS:1001457)
# K Data File Name
$ws.Range("E1").Value = "Parameter [analysis software]"
$ws.Range("F1").Value = "Term Source REF (MS:1001456)"
$ws.Range("G1").Value = "Term Accession Number (MS:1001456)"
$ws.Range("H1").Value = "Parameter [data processing software]"
$ws.Range("I1").Value = "Term Source REF (MS:1001457)"
$ws.Range("J1").Value = "Term Accession Number (MS:1001457)"
$ws.Range("K1").Value = "Data File Name"

# Fill content examples taken from older templates into the three
# "software" columns (acquisition / analysis / data processing).
$acquisitionSoftware = @(
    "SCIEX TOF/TOF Series Explorer Software",
    "6300 Series Ion Trap Data Analysis Software",
    "Analyst",
    "ChromaTOF HRT software",
    "Xcalibur",
    "apexControl"
)
$analysisSoftware = @(
    "Analyst",
    "ChromaTOF software",
    "Comet",
    "Compass",
    "BioAnalyst",
    "BioTools"
)
$dataProcessingSoftware = @(
    "Analyst",
    "BioAnalyst",
    "BioTools",
    "Bioworks",
    "CAMERA",
    "Data Explorer"
)

for ($i = 0; $i -lt 6; $i++) {
    $r = 2 + $i
    $ws.Range("B$r").Value = $acquisitionSoftware[$i]
    $ws.Range("E$r").Value = $analysisSoftware[$i]
    $ws.Range("H$r").Value = $dataProcessingSoftware[$i]
}

# Column widths / visibility to mirror the other ontology-term columns:
# the Term Source REF / Term Accession Number columns stay hidden, the
# "Parameter [...]" columns stay visible with a best-fit width.
# Columns 6 (F) and 7 (G) already carry the correct hidden REF/Accession
# widths inherited from the former columns 6/7, so they are left alone.
$ws.Columns.Item(2).ColumnWidth = 31.877604166666668   # -> xml width 32.7109375

$ws.Columns.Item(5).ColumnWidth = 29.307291666666668   # -> xml width 30.140625

$ws.Columns.Item(8).ColumnWidth = 36.166666666666664   # -> xml width 37

$ws.Columns.Item(9).Hidden = $true
$ws.Columns.Item(9).ColumnWidth = -0.8333333333333334  # -> xml width 0

$ws.Columns.Item(10).Hidden = $true
$ws.Columns.Item(10).ColumnWidth = -0.8333333333333334 # -> xml width 0

$ws.Columns.Item(11).ColumnWidth = 15.877604166666666  # -> xml width 16.7109375

# ---------------------------------------------------------------------
# 2) Template metadata sheet: bump version 1.1.3 -> 1.1.4
# ---------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("SwateTemplateMetadata")
$wsMeta.Range("B3").Value = "1.1.4"
